$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Semestre ideal" value from "EF-7,EM-8" to "EF-7"
$ws.Range("B9").Value = "EF-7"
$ws.Range("C9").Value = "EF-7"

# Remove the "Requisitos:" rows (23 and 24), which held the
# course-requirement text that's no longer part of the sheet.
$ws.Rows("23:24").Delete()
